# Scheduled-runner style refresh of the Leve-flipping profit data across
# the crafting-job sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). Updates the
# market-price / profit columns (H:N) for the rows whose prices moved
# since the last run. A handful of rows gain or lose their profit cells
# (M/N) entirely depending on whether NQ/HQ pricing was available this
# pass, so those are cleared explicitly rather than zeroed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 901
$ws.Range("I43").Value = 901
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 901
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -832
$ws.Range("N43").ClearContents()
$ws.Range("H88").Value = 25799.4
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 25799.4
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 25799.4
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -26611.4
$ws.Range("H91").Value = 25799.4
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 25799.4
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 25799.4
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -28607.4
$ws.Range("H113").Value = 25002282
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 31252102
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 31252102
$ws.Range("M113").Value = 254
$ws.Range("N113").Value = -31258610
$ws.Range("H140").Value = 65783.336
$ws.Range("J140").Value = 65783.336
$ws.Range("L140").Value = 65783.336
$ws.Range("N140").Value = -76143.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15863.903
$ws.Range("I32").Value = 11759.274
$ws.Range("J32").Value = 24073.162
$ws.Range("K32").Value = 11759.274
$ws.Range("L32").Value = 24073.162
$ws.Range("M32").Value = -11472.274
$ws.Range("N32").Value = -24647.162
$ws.Range("H61").Value = 347949
$ws.Range("I61").Value = 2716.4211
$ws.Range("J61").Value = 1003890.9
$ws.Range("K61").Value = 2716.4211
$ws.Range("L61").Value = 1003890.9
$ws.Range("M61").Value = -2504.4211
$ws.Range("N61").Value = -1004314.9
$ws.Range("H88").Value = 2406.3333
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 2609.5
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 2609.5
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -3421.5
$ws.Range("H91").Value = 2406.3333
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 2609.5
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 2609.5
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -5417.5
$ws.Range("H110").Value = 7025.5537
$ws.Range("I110").Value = 8083.7827
$ws.Range("J110").Value = 2157.7
$ws.Range("K110").Value = 8083.7827
$ws.Range("L110").Value = 2157.7
$ws.Range("M110").Value = -6038.7827
$ws.Range("N110").Value = -6247.7
$ws.Range("H122").Value = 2568983.8
$ws.Range("I122").Value = 2853870.8
$ws.Range("K122").Value = 8561612.399999999
$ws.Range("M122").Value = -8559162.399999999
$ws.Range("H136").Value = 347949
$ws.Range("I136").Value = 2716.4211
$ws.Range("J136").Value = 1003890.9
$ws.Range("K136").Value = 8149.263300000001
$ws.Range("L136").Value = 3011672.7
$ws.Range("M136").Value = -5599.263300000001
$ws.Range("N136").Value = -3016772.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2127.25
$ws.Range("I20").Value = 2069.6667
$ws.Range("J20").Value = 2300
$ws.Range("K20").Value = 2069.6667
$ws.Range("L20").Value = 2300
$ws.Range("M20").Value = -1822.6667
$ws.Range("N20").Value = -2794
$ws.Range("H107").Value = 375619
$ws.Range("I107").Value = 375619
$ws.Range("K107").Value = 375619
$ws.Range("M107").Value = -373699
$ws.Range("H134").Value = 38130.035
$ws.Range("I134").Value = 2524.1155
$ws.Range("K134").Value = 7572.3465
$ws.Range("M134").Value = -5037.3465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9264192
$ws.Range("I31").Value = 1597.963
$ws.Range("J31").Value = 18526786
$ws.Range("K31").Value = 1597.963
$ws.Range("L31").Value = 18526786
$ws.Range("M31").Value = -1302.963
$ws.Range("N31").Value = -18527376
$ws.Range("H34").Value = 9264192
$ws.Range("I34").Value = 1597.963
$ws.Range("J34").Value = 18526786
$ws.Range("K34").Value = 1597.963
$ws.Range("L34").Value = 18526786
$ws.Range("M34").Value = -1395.963
$ws.Range("N34").Value = -18527190
$ws.Range("H134").Value = 10732143
$ws.Range("I134").Value = 15155406
$ws.Range("J134").Value = 1000963.9
$ws.Range("K134").Value = 45466218
$ws.Range("L134").Value = 3002891.7
$ws.Range("M134").Value = -45463683
$ws.Range("N134").Value = -3007961.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1725125.6
$ws.Range("I131").Value = 5883011
$ws.Range("J131").Value = 1124.3903
$ws.Range("K131").Value = 17649033
$ws.Range("L131").Value = 3373.1709
$ws.Range("M131").Value = -17643993
$ws.Range("N131").Value = -13453.1709
$ws.Range("H137").Value = 17676.791
$ws.Range("I137").Value = 9854.286
$ws.Range("J137").Value = 28628.3
$ws.Range("K137").Value = 29562.858
$ws.Range("L137").Value = 85884.89999999999
$ws.Range("M137").Value = -24462.858
$ws.Range("N137").Value = -96084.89999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2720.6316
$ws.Range("I102").Value = 2192.6924
$ws.Range("J102").Value = 3864.5
$ws.Range("K102").Value = 2192.6924
$ws.Range("L102").Value = 3864.5
$ws.Range("M102").Value = -570.6923999999999
$ws.Range("N102").Value = -7108.5
$ws.Range("H123").Value = 28997
$ws.Range("J123").Value = 28997
$ws.Range("L123").Value = 28997
$ws.Range("N123").Value = -33897

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2471.8823
$ws.Range("I7").Value = 1562.4615
$ws.Range("K7").Value = 1562.4615
$ws.Range("M7").Value = -1450.4615
$ws.Range("H40").Value = 2681.818
$ws.Range("I40").Value = 2611.111
$ws.Range("K40").Value = 2611.111
$ws.Range("M40").Value = -2475.111
$ws.Range("H126").Value = 2471.8823
$ws.Range("I126").Value = 1562.4615
$ws.Range("K126").Value = 4687.3845
$ws.Range("M126").Value = -2217.3845
$ws.Range("H132").Value = 3626.9714
$ws.Range("I132").Value = 3468.5
$ws.Range("J132").Value = 4577.8
$ws.Range("K132").Value = 10405.5
$ws.Range("L132").Value = 13733.4
$ws.Range("M132").Value = -7875.5
$ws.Range("N132").Value = -18793.4
$ws.Range("H138").Value = 77651.14
$ws.Range("J138").Value = 77651.14
$ws.Range("L138").Value = 77651.14
$ws.Range("N138").Value = -87931.14

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 43478956
$ws.Range("I107").Value = 62500604
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 187501812
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = -187499892
$ws.Range("N107").Value = -6540
$ws.Range("H122").Value = 1155.4445
$ws.Range("I122").Value = 1090.8
$ws.Range("J122").Value = 1236.25
$ws.Range("K122").Value = 3272.4
$ws.Range("L122").Value = 3708.75
$ws.Range("M122").Value = -822.3999999999996
$ws.Range("N122").Value = -8608.75
